$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openml_100")

# The "STC" study column (column E, header "STC") is being removed entirely.
# Deleting the whole column shifts F:L left into E:K, and removes the now
# unused "STC" entry from the shared string table.
$ws.Columns.Item(5).Delete()

# Update the view: active selection moves to D2 and normal zoom is pinned to 100.
$ws.Activate()
$ws.Range("D2").Select()
$win = $excel.ActiveWindow
$win.Zoom = 100
